# Auto-generated edit script: updates cryptos price/volume table per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.897.78'
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("D3").Value = '1.874.83'
$ws.Range("D4").Value = '''0.9996'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''0.7378'
$ws.Range("D6").Value = '''242.09'
$ws.Range("E6").Value = '  -0.75%  '
$ws.Range("D7").Value = '''0.9989'
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '''0.3156'
$ws.Range("E8").Value = '  +0.85%  '
$ws.Range("D9").Value = '''0.07181'
$ws.Range("E9").Value = '  -0.96%  '
$ws.Range("D10").Value = '''24.72'
$ws.Range("E10").Value = '  -4.40%  '
$ws.Range("D11").Value = '''0.08360'
$ws.Range("E11").Value = '  -4.01%  '
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = '''0.7501'
$ws.Range("E12").Value = '  -3.08%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''5.419'
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.879.56'
$ws.Range("E14").Value = '  -11.54%  '
$ws.Range("D15").Value = '''92.52'
$ws.Range("E15").Value = '  -2.16%  '
$ws.Range("D16").Value = '29.893.88'
$ws.Range("E16").Value = '  -1.52%  '
$ws.Range("D17").Value = '''6.066'
$ws.Range("D18").Value = '''246.11'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("E19").Value = '  -2.52%  '
$ws.Range("D20").Value = '''0.000007830'
$ws.Range("E20").Value = '  -0.41%  '
$ws.Range("D21").Value = '''0.9985'
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("D22").Value = '2.127.61'
$ws.Range("E22").Value = '  -11.34%  '
$ws.Range("D23").Value = '''7.986'
$ws.Range("E23").Value = '  -1.69%  '
$ws.Range("D24").Value = '''1.000'
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = '''0.1551'
$ws.Range("E25").Value = '  -7.01%  '
$ws.Range("D26").Value = '''9.260'
$ws.Range("D27").Value = '''165.08'
$ws.Range("E27").Value = '  +1.02%  '
$ws.Range("D28").Value = '''18.65'
$ws.Range("E28").Value = '  -0.99%  '
$ws.Range("D29").Value = '''2.031'
$ws.Range("E29").Value = '  -1.13%  '
$ws.Range("D30").Value = '''1.508'
$ws.Range("E30").Value = '  +5.14%  '
$ws.Range("D31").Value = '''4.586'
$ws.Range("E31").Value = '  +1.62%  '
$ws.Range("D32").Value = '''1.534'
$ws.Range("E32").Value = '  -0.57%  '
$ws.Range("E33").Value = '  +3.15%  '
$ws.Range("D34").Value = '''0.05318'
$ws.Range("E34").Value = '  -2.87%  '
$ws.Range("D35").Value = '''1.238'
$ws.Range("E35").Value = '  -0.57%  '
$ws.Range("D36").Value = '''0.7545'
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").Value = '''1.000'
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("E39").Value = '  -0.72%  '
$ws.Range("D40").Value = '''2.753'
$ws.Range("E40").Value = '  -1.21%  '
$ws.Range("D41").Value = '''0.4513'
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '1.113.04'
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").Value = '''6.046'
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("D44").Value = '''72.36'
$ws.Range("E44").Value = '  -1.70%  '
$ws.Range("D45").Value = '''0.8527'
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").Value = '''103.23'
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("D48").Value = '''7.640'
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").Value = '''3.093'
$ws.Range("E49").Value = '  +2.78%  '
$ws.Range("D50").Value = '''1.840'
$ws.Range("D51").Value = '2.023.43'
$ws.Range("E51").Value = '  -9.79%  '
